$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 31
$lastColBefore = 16   # column P
$lastColAfter = 17    # column Q

# --- Row 1 (header row): a new column is inserted at column B (index 2). ---
# Shift existing header cells B..P (2..16) right by one, into C..Q (3..17).
for ($c = $lastColAfter; $c -ge 3; $c--) {
    $src = $ws.Cells.Item(1, $c - 1)
    $dst = $ws.Cells.Item(1, $c)
    $dst.Value = $src.Value()
}
# New header cell at B1.
$ws.Cells.Item(1, 2).Value = "Unnamed: 0.3"

# --- Data rows 2..31: a new column is inserted at column E (index 5). ---
# Shift existing cells E..P (5..16) right by one, into F..Q (6..17).
for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = $lastColAfter; $c -ge 6; $c--) {
        $src = $ws.Cells.Item($r, $c - 1)
        $dst = $ws.Cells.Item($r, $c)
        $dst.Value = $src.Value()
    }
    # New cell at column E duplicates the (already-shifted-in-place) column D value.
    $dCell = $ws.Cells.Item($r, 4)
    $ws.Cells.Item($r, 5).Value = $dCell.Value()
}
